# Swap the order of the "System" / email entries in column G:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# This applies to every row in the used range where column G currently
# holds that exact text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
